$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Insert a new row for the "Blur" command right after the "Clear" row (row 9),
# pushing every row from the old row 10 onward down by one. This reproduces
# the row-shift seen across the whole sheet in the diff.
# ---------------------------------------------------------------------------
$ws.Rows(10).Insert()

# Match the formatting used by the surrounding data rows ("Good" cell style,
# internally style index 1) for the newly inserted cells only (A:F, since
# that's all the new row uses).
$ws.Range("A10:F10").Style = "Good"

$ws.Range("A10").Value = "Blur"
$ws.Range("B10").Value = "Character/Background/All"
$ws.Range("C10").Value = "gradual? Relative?"
$ws.Range("D10").Value = "start_blur"
$ws.Range("E10").Value = "end_blur"
$ws.Range("F10").Value = "time"

# ---------------------------------------------------------------------------
# "allowed multiple flags to be triggered for each choice": the flag1 column
# on the Display/Choice row (shifted from row 17 to row 18 by the insert
# above) gets a clarified label.
# ---------------------------------------------------------------------------
$ws.Range("F18").Value = "flag1 (use comma (,) to separate mutliple flags)"

# ---------------------------------------------------------------------------
# View-state tweaks captured in the diff: zoomed in a bit, and the selection
# left on F19 instead of the previous scroll position / selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
$ws.Range("F19").Select()
